$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each 6-row block (one per day-group), a new daily reading is prepended:
# the old block-start row is discarded, the old 2nd row becomes the new 3rd row,
# the old 3rd row becomes the new 4th row (old 4th row is discarded), and the
# first two rows get the new readings (updated "consum mitja diari" calc).
$blockStarts = @(1, 7, 13, 19, 25, 31, 37, 43, 49, 55)

foreach ($r0 in $blockStarts) {
    $r1 = $r0 + 1
    $r2 = $r0 + 2
    $r3 = $r0 + 3

    # capture the existing values of rows 2 and 3 of the block before overwriting
    $oldA1 = $ws.Cells.Item($r1, 1).Value2
    $oldB1 = $ws.Cells.Item($r1, 2).Value2
    $oldC1 = $ws.Cells.Item($r1, 3).Value2
    $oldA2 = $ws.Cells.Item($r2, 1).Value2
    $oldB2 = $ws.Cells.Item($r2, 2).Value2
    $oldC2 = $ws.Cells.Item($r2, 3).Value2

    # shift down one row: old row1 -> row2 (row3 in sheet), old row2 -> row3 (row4 in sheet)
    $ws.Cells.Item($r2, 1).Value2 = $oldA1
    $ws.Cells.Item($r2, 2).Value2 = $oldB1
    $ws.Cells.Item($r2, 3).Value2 = $oldC1
    $ws.Cells.Item($r3, 1).Value2 = $oldA2
    $ws.Cells.Item($r3, 2).Value2 = $oldB2
    $ws.Cells.Item($r3, 3).Value2 = $oldC2
}

# new values for the top two rows of each block
$newTop = @{
    1 = @(@(9235, 7162, 7077), @(6668, 2856, 2973))
    7 = @(@(8467, 6405, 6598), @(5513, 2471, 2556))
    13 = @(@(9386, 6289, 6385), @(5370, 2316, 2309))
    19 = @(@(7933, 4609, 4624), @(4364, 1440, 1579))
    25 = @(@(10566, 7921, 8019), @(5873, 2606, 2656))
    31 = @(@(14563, 9290, 9047), @(6580, 3214, 3375))
    37 = @(@(22846, 13406, 11278), @(8977, 3985, 4141))
    43 = @(@(13791, 13088, 11657), @(8247, 4344, 4509))
    49 = @(@(9795, 9224, 8784), @(7327, 3234, 3381))
    55 = @(@(9961, 8720, 8305), @(5831, 2479, 2583))
}

foreach ($r0 in $blockStarts) {
    $r1 = $r0 + 1
    $pair = $newTop[$r0]
    $row0 = $pair[0]
    $row1 = $pair[1]
    $ws.Cells.Item($r0, 1).Value2 = $row0[0]
    $ws.Cells.Item($r0, 2).Value2 = $row0[1]
    $ws.Cells.Item($r0, 3).Value2 = $row0[2]
    $ws.Cells.Item($r1, 1).Value2 = $row1[0]
    $ws.Cells.Item($r1, 2).Value2 = $row1[1]
    $ws.Cells.Item($r1, 3).Value2 = $row1[2]
}
